$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the period after "5.6.1.1" in the Russian title (B1)
$ws.Range("B1").Value = "5.6.1.1 Доля замужних женщин и сексуально активных не замужних женщин в возрасте 15-49 лет, которые были осведомлены о соврменном методе контрацепции"

# Update the "urban" row (row 6) with fuller wording
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

# Update the "rural" row (row 7) with fuller wording
$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Update the selected range shown when the workbook is reopened
$ws.Range("A6:C7").Select()
